# "Update presentation. Try out 3d plot"
# Swap the subject-ID assignments (column A) between the two K-means cluster
# sheets for rows 2-208 (the rows present in both sheets). Cluster1.K keeps
# its extra rows 209-261 untouched; column B (the cluster number) is left
# alone in both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Cluster1.K")
$ws2 = $wb.Worksheets.Item("Cluster2.K")

$firstRow = 2
$lastRow = 208

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell1 = $ws1.Cells.Item($r, 1)
    $cell2 = $ws2.Cells.Item($r, 1)

    $val1 = $cell1.Value()
    $val2 = $cell2.Value()

    $cell1.Value = $val2
    $cell2.Value = $val1
}
